# Generate Report for Handoff
#
# The localization-status report is regenerated: the previous run had
# handed the files back ("Handed back: in sync with en-US") and is now
# ready to be hand off again ("Ready for handoff"). The handoff
# timestamps on the Overview sheet and on the per-locale (zh-cn / de-de)
# detail sheets are refreshed to the new generation time, and the
# "Status" columns (which now hold shorter text) are resized to fit.

$wb = $excel.ActiveWorkbook

$statusOld = "Handed back: in sync with en-US"
$statusNew = "Ready for handoff"

# ---- Overview sheet ---------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

# zh-cn / de-de status columns
$ov.Range("E2").Value = $statusNew
$ov.Range("F2").Value = $statusNew

# Latest HO Xliff Generate Date
$ov.Range("G2").Value = "2016-08-19 06:57:45"

# Status columns got narrower now that the text is shorter
$ov.Columns.Item(5).ColumnWidth = 16.3
$ov.Columns.Item(6).ColumnWidth = 16.3

# ---- zh-cn detail sheet ------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("C2").Value = $statusNew
$zh.Range("H2").Value = "2016-08-19 06:57:40"
$zh.Columns.Item(3).ColumnWidth = 16.3

# ---- de-de detail sheet ------------------------------------------------
$de = $wb.Worksheets.Item("de-de")
$de.Range("C2").Value = $statusNew
$de.Range("H2").Value = "2016-08-19 06:57:45"
$de.Columns.Item(3).ColumnWidth = 16.3
